$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new PR log entry (row 11) from PR #33
$ws.Range("A11").Value = 33
$ws.Range("B11").Value = "Edit2"
$ws.Range("C11").Value = "riya-morankar"
$ws.Range("D11").Value = "N/A"
$ws.Range("E11").Value = "edit2 to main"

# F11 ("2025-06-18") looks like a date, and Excel would normally parse it
# into a date serial. Force it to be entered as literal text by
# pre-formatting the cell, then strip the format override afterwards so the
# cell keeps the plain "General" styling used by the rest of the sheet.
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "2025-06-18"
$ws.Range("F11").ClearFormats()
